# Update "想去人数" (number of interested people) counts for several
# convention entries across sheets, reflecting the latest scrape output.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 700
$ws1.Range("F12").Value = 32985
$ws1.Range("F13").Value = 7067
$ws1.Range("F15").Value = 366
$ws1.Range("F16").Value = 573
$ws1.Range("F24").Value = 801
$ws1.Range("F28").Value = 442

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 1189

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 1189
$ws4.Range("F9").Value  = 700
$ws4.Range("F21").Value = 7067
$ws4.Range("F23").Value = 366
$ws4.Range("F25").Value = 573
$ws4.Range("F33").Value = 801
$ws4.Range("F36").Value = 442
